$d = $word.ActiveDocument

$r = $d.Content
$r.Find.Execute("349÷9=38, 7", $false, $false, $false, $false, $false, $true, 1, $false, "480÷9=53, 3", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("422÷9=46, 8", $false, $false, $false, $false, $false, $true, 1, $false, "435÷2=217, 1", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("715÷5=143, 0", $false, $false, $false, $false, $false, $true, 1, $false, "664÷6=110, 4", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("196÷9=21, 7", $false, $false, $false, $false, $false, $true, 1, $false, "314÷7=44, 6", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("649÷4=162, 1", $false, $false, $false, $false, $false, $true, 1, $false, "295÷5=59, 0", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("304÷2=152, 0", $false, $false, $false, $false, $false, $true, 1, $false, "857÷3=285, 2", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("580÷7=82, 6", $false, $false, $false, $false, $false, $true, 1, $false, "939÷3=313, 0", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("115÷7=16, 3", $false, $false, $false, $false, $false, $true, 1, $false, "480÷8=60, 0", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("777÷8=97, 1", $false, $false, $false, $false, $false, $true, 1, $false, "715÷5=143, 0", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("560÷5=112, 0", $false, $false, $false, $false, $false, $true, 1, $false, "703÷9=78, 1", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("572÷6=95, 2", $false, $false, $false, $false, $false, $true, 1, $false, "545÷7=77, 6", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("931÷8=116, 3", $false, $false, $false, $false, $false, $true, 1, $false, "763÷4=190, 3", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("784÷8=98, 0", $false, $false, $false, $false, $false, $true, 1, $false, "491÷3=163, 2", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("956÷2=478, 0", $false, $false, $false, $false, $false, $true, 1, $false, "848÷3=282, 2", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("191÷5=38, 1", $false, $false, $false, $false, $false, $true, 1, $false, "562÷8=70, 2", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("690÷8=86, 2", $false, $false, $false, $false, $false, $true, 1, $false, "588÷8=73, 4", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("864÷7=123, 3", $false, $false, $false, $false, $false, $true, 1, $false, "954÷4=238, 2", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("777÷2=388, 1", $false, $false, $false, $false, $false, $true, 1, $false, "471÷6=78, 3", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("390÷8=48, 6", $false, $false, $false, $false, $false, $true, 1, $false, "983÷4=245, 3", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("476÷3=158, 2", $false, $false, $false, $false, $false, $true, 1, $false, "440÷4=110, 0", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("659÷7=94, 1", $false, $false, $false, $false, $false, $true, 1, $false, "518÷5=103, 3", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("267÷8=33, 3", $false, $false, $false, $false, $false, $true, 1, $false, "421÷3=140, 1", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("577÷7=82, 3", $false, $false, $false, $false, $false, $true, 1, $false, "864÷3=288, 0", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("535÷3=178, 1", $false, $false, $false, $false, $false, $true, 1, $false, "879÷6=146, 3", 1) | Out-Null
$r = $d.Content
$r.Find.Execute("962÷6=160, 2", $false, $false, $false, $false, $false, $true, 1, $false, "860÷6=143, 2", 1) | Out-Null
